# Generate Report for Handoff
# Updates the localization-status report to reflect that the
# f674905b-9063-4abe-af0b-b694f7aa8537.md file is now ready for handoff
# (latest handback file was stale), across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2445d03f9f357a8dfff44744c36a2b239c26b522/e2e/f674905b-9063-4abe-af0b-b694f7aa8537.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5223164ee02894205d0aae81228dac05499d33b/e2e/f674905b-9063-4abe-af0b-b694f7aa8537.md."

# --- Overview sheet: row 3 is the f674905b-...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = "2016-08-21 08:53:44"

# --- zh-cn sheet: row 3 is the f674905b-...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("H3").Value = "2016-08-21 08:53:40"
$wsZhCn.Range("P3").Value = $errorDetail
# 39.17 chars round-trips (via the engine's px-based column-width model) to a
# stored OOXML width of exactly 40, matching columns A/G/I/J on this sheet.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the f674905b-...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("H3").Value = "2016-08-21 08:53:44"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17

Write-Output "Report updated for handoff"
